# Timesheet.xlsx edit: "proofreading and smaller changes to chapter around CP"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: updated hours worked (column D) for the CP chapter weeks.
# ---------------------------------------------------------------------------
$ws1.Range("D13").Value = 30
$ws1.Range("D14").Value = 40
$ws1.Range("D15").Value = 40
$ws1.Range("D17").Value = 20
$ws1.Range("D18").Value = 40
$ws1.Range("D19").Value = 40
$ws1.Range("D20").Value = 40
$ws1.Range("D21").Value = 40
$ws1.Range("D22").Value = 40
$ws1.Range("D23").Value = 40
$ws1.Range("D24").Value = 40
$ws1.Range("D25").Value = 40
$ws1.Range("D26").Value = 40
$ws1.Range("D27").Value = 40
$ws1.Range("D28").Value = 40
$ws1.Range("D29").Value = 40

# Row 13's "status" cell (H13) flips from "tbd" to "proofreading", and picks
# up the same bottom-border style ("s=4") as the rest of that (bottom-border)
# row.
$ws1.Range("G13").Copy()
$ws1.Range("H13").PasteSpecial(-4122) # xlPasteFormats
$ws1.Range("H13").Value = "proofreading"

# B34 now compares against F13 (the CP-chapter checkpoint row) instead of F12.
$ws1.Range("B34").Formula = '=$G$31-F13'

# ---------------------------------------------------------------------------
# Sheet2: extra logged hours + note for week 38 (row 12), and C11 gains some
# trailing zero terms (no value change).
# ---------------------------------------------------------------------------
$ws2.Range("C11").Formula = "=8+7+5+6+0+0+0"
$ws2.Range("C12").Formula = "=0+7"
$ws2.Range("G12").ClearContents()
$ws2.Range("E12").Value = "Writing conclusion of chapter CP, proofreading chapter CP, partial reading 73 Simplify A Theorem Prover for Program Checking 1066100.1066102.pdf"

# ---------------------------------------------------------------------------
# Move/resize the progress chart on Sheet1 (it was nudged right slightly).
# ---------------------------------------------------------------------------
$co = $ws1.ChartObjects().Item(1)
$co.Left = 604.1044921875
$co.Width = 433.0625

# ---------------------------------------------------------------------------
# View state: Sheet2 becomes the active/selected sheet, with C12 selected;
# Sheet1 keeps a selection of B35 for when it's revisited.
# ---------------------------------------------------------------------------
$ws1.Range("B35").Select()
$ws2.Range("C12").Select()
$ws2.Activate()
